# Apply updated crypto price/volume data per Thu Aug 15 17:40:35 UTC 2024 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.201.89"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "2.635.29"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'529.67"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("D6").Value = "'144.87"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("D9").Value = "'6.64"
$ws.Range("E9").Value = "  -5.80%  "
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("E12").Value = "  +0.46%  "
$ws.Range("D13").Value = "3.099.72"
$ws.Range("E13").Value = "  -0.70%  "
$ws.Range("D14").Value = "59.180.16"
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.726.05"
$ws.Range("E15").Value = "  +3.16%  "
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Value = "'20.60"
$ws.Range("E16").Value = "  -2.92%  "
$ws.Range("D17").Value = "'0.0000135"
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("D18").Value = "'342.49"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").Value = "'4.40"
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").Value = "'10.49"
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "'66.53"
$ws.Range("E23").Value = "  +3.94%  "
$ws.Range("D24").Value = "'0.412"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").Value = "2.758.38"
$ws.Range("E26").Value = "  -0.61%  "
$ws.Range("D27").Value = "'0.998"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").Value = "0.0₃0787"
$ws.Range("E29").Value = "  -2.13%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").Value = "'6.28"
$ws.Range("E31").Value = "  -6.06%  "
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("D33").Value = "'18.92"
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("D34").Value = "'149.98"
$ws.Range("E34").Value = "  +0.58%  "
$ws.Range("E35").Value = "  -1.32%  "
$ws.Range("E36").Value = "  -3.16%  "
$ws.Range("E37").Value = "  -0.89%  "
$ws.Range("D38").Value = "'0.837"
$ws.Range("E38").Value = "  -7.40%  "
$ws.Range("D39").Value = "'0.837"
$ws.Range("E39").Value = "  -4.87%  "
$ws.Range("E40").Value = "  -2.22%  "
$ws.Range("D41").Value = "'3.59"
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("D42").Value = "'0.997"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").Value = "'0.0973"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").Value = "'0.596"
$ws.Range("E44").Value = "  -3.67%  "
$ws.Range("D46").Value = "'268.27"
$ws.Range("E46").Value = "  -2.63%  "
$ws.Range("D47").Value = "'19.03"
$ws.Range("E47").Value = "  -4.24%  "
$ws.Range("D48").Value = "'0.0530"
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("D49").Value = "2.028.35"
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "'4.68"
$ws.Range("E50").Value = "  -1.81%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0227"
$ws.Range("E51").Value = "  -0.59%  "
